# Weekly data refresh: insert the latest week's "Arándano (blue)" price
# record at the top of the time series (row 11) and push the previously
# recorded weeks down by one row, matching the new-week-on-top layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the new weekly record by inserting a row at row 11
# (the first data row after the existing "Primera" entry for 2021-12-07
# that is about to move down).
$ws.Rows("11:11").Insert()

# Duplicate the row that was just pushed down to row 12 into the new
# row 11 so every static column (Mercado, Región, Producto, Calidad,
# Unidad de comercialización, Origen, Kg/unidad, styles, ...) keeps the
# same values/formatting as the series it belongs to.
$ws.Range("A12:T12").Copy($ws.Range("A11:T11"))

# Now overwrite the fields that differ for this new week's record.
$ws.Cells.Item(11, 4).Value = 44592   # Fecha
$ws.Cells.Item(11, 13).Value = 150    # Volumen
$ws.Cells.Item(11, 14).Value = 3500   # Precio mínimo
$ws.Cells.Item(11, 15).Value = 3500   # Precio máximo
$ws.Cells.Item(11, 16).Value = 3500   # Precio promedio ponderado
$ws.Cells.Item(11, 19).Value = 1750   # Precio $/Kg
